$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "42.615.33"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.02%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "2.526.92"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.71%  "

$ws.Range("E4").Value = "  -0.04%  "

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "315.06"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.00%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "98.41"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -1.33%  "

$ws.Range("E7").Value = "  -1.25%  "

$ws.Range("E8").Value = "  +0.00%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.519"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.37%  "

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "35.19"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.01%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.0800"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.39%  "

$ws.Range("E12").Value = "  +1.44%  "

$ws.Range("E13").Value = "  -1.85%  "

$ws.Range("E14").Value = "  -1.04%  "

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "2.545.04"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.50%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "15.21"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -6.05%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.811"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -3.42%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "42.619.03"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.02%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.59"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -2.57%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.0₃0939"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -1.11%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "12.09"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.68%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "69.14"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.17%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "241.97"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.15%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "2.86"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -1.24%  "

$ws.Range("E25").Value = "  -3.21%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.08%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "25.53"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -3.05%  "

$ws.Range("E28").Value = "  -4.28%  "

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "10.01"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.61%  "

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "37.65"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -5.75%  "

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "5.92"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +4.44%  "

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "155.80"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -1.70%  "

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "2.71"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -2.76%  "

$ws.Range("E34").Value = "  +0.93%  "

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.0782"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.89%  "

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "3.13"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.16%  "

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.97"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -3.62%  "

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "17.55"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.18%  "

$ws.Range("E39").Value = "  -2.46%  "

$ws.Range("E40").Value = "  -0.76%  "

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "4.24"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +0.24%  "

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "21.69"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.62%  "

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "2.025.65"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +3.36%  "

$ws.Range("E45").Value = "  -0.39%  "

$ws.Range("E46").Value = "  -3.89%  "

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "8.86"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.60%  "

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.768.38"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -1.06%  "

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "80.16"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.73%  "

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "71.93"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -0.54%  "

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.188"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -2.43%  "
